$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Capture the existing B/C values (rows 6..20, 1-indexed) before we
    # overwrite anything, since the new B2:C16 block is just the old
    # B6:C20 block shifted up by 4 rows. Column A (the "Cutoff" index
    # column) stays exactly as-is for rows 2..16.
    $oldB = @{}
    $oldC = @{}
    for ($r = 6; $r -le 20; $r++) {
        $oldB[$r] = $ws.Cells.Item($r, 2).Value2
        $oldC[$r] = $ws.Cells.Item($r, 3).Value2
    }

    # Write the shifted B/C values into rows 2..16.
    for ($newR = 2; $newR -le 16; $newR++) {
        $srcR = $newR + 4
        $ws.Cells.Item($newR, 2).Value = $oldB[$srcR]
        $ws.Cells.Item($newR, 3).Value = $oldC[$srcR]
    }

    # Rows 17..20 are no longer part of the table; clear them so the
    # worksheet's used range shrinks back down to A1:C16.
    $ws.Range("A17:C20").Clear()
}
